# Add a new risk row to the end of the risk-assessment table:
# "Files in wrong directory" / ... / medium / medium / medium
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "Files in wrong directory"
$newRow.Cells.Item(2).Range.Text = "If I accidentally save an important file for the project in the wrong file it will not be pushed to my github repo"
$newRow.Cells.Item(3).Range.Text = "Try to save things in the correct place first time, every week check contents of project folder against specification checklist"
$newRow.Cells.Item(4).Range.Text = "Ensure that all files are in the correct place so that they will be pushed correctly."
$newRow.Cells.Item(5).Range.Text = "medium"
$newRow.Cells.Item(6).Range.Text = "medium"
$newRow.Cells.Item(7).Range.Text = "medium"
